$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: AS4 v1.0 is now marked Deprecated (E4 = TRUE) with a
#     "Deprecated since?" value of 3 (F4, a literal number) ---
$ws.Range("E4").Formula = "=TRUE"
$ws.Range("F4").Value = 3

# --- Row 5 (new): AS4 v2.0 identifier ---
$ws.Range("A5").Value = "AS4"
$ws.Range("C5").Value = "peppol-transport-as4-v2_0"
$ws.Range("E5").Value = $false

# B5 ("2.0") and D5 ("3") must stay as TEXT even though they look like
# numbers, matching the sheet's existing text-typed version/since columns.
# Temporarily force text format, enter the value, then restore the
# general-format look of column B by re-applying A5's format (column D
# keeps its own text-format default, so no further fix-up is needed there).
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "2.0"
$ws.Range("A5").Copy()
$ws.Range("B5").PasteSpecial(-4122)

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "3"

# E5 must carry the same style as the other "Deprecated" cells (E2:E4).
$ws.Range("E4").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value = $false

# Match the author's final selection (A4:B4) recorded in the saved file.
$ws.Range("A4:B4").Select() | Out-Null
